$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear row 3 (Item2) - we will rewrite A1:C2 and remove row 3 content
$ws.Rows.Item(3).ClearContents()

# Row 1: numeric values and formula
$ws.Range("A1").Value = 10
$ws.Range("B1").Value = 20
$ws.Range("C1").Formula = "=A1+B1"

# Row 2: numeric values and formula
$ws.Range("A2").Value = 5
$ws.Range("B2").Value = 15
$ws.Range("C2").Formula = "=SUM(A2:B2)"
